# Update column B (period/date reference) text on the "dataText" sheet:
# remove the trailing "North East and North of Tyne MCA have temporarily
# been removed, on account of North of Tyne being integrated into North
# East from now on" sentence, leaving just the data-period prefix.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataText")

$values = @{
    2  = "Oct 2023-Sep 2024. "
    3  = "Oct 2023-Sep 2024. "
    4  = "Oct 2023-Sep 2024. "
    5  = "Oct 2023-Sep 2024. "
    6  = "Oct 2023-Sep 2024. "
    7  = "Oct 2023-Sep 2024. "
    8  = "Oct 2023-Sep 2024. "
    9  = "Oct 2023-Sep 2024. "
    10 = "Nov 2024 data."
    11 = "Mar 2024 data. "
    12 = "Dec 2022 - Dec 2023 data. "
    13 = "Dec 2022 - Dec 2023 data. "
    14 = "AY23/24 data. "
    15 = "AY23/24 data. "
    16 = "AY23/24 data. "
    17 = "AY23/24 data. "
    18 = "Jan-Dec 2023 data. "
    19 = "Jan-Dec 2023 data. "
    20 = "AY22/23 data. "
    21 = "AY22/23 data. "
    22 = "Growth from 2023 to 2035. "
    23 = "AY22/23 data. "
    24 = "AY22/23 data. "
}

foreach ($row in $values.Keys) {
    $ws.Range("B$row").Value = $values[$row]
}
